$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1870.6522
$ws.Range("J40").Value = 1744.2
$ws.Range("L40").Value = 1744.2
$ws.Range("N40").Value = -2094.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3087
$ws.Range("I116").Value = 1491.6666
$ws.Range("J116").Value = 5480
$ws.Range("K116").Value = 1491.6666
$ws.Range("L116").Value = 5480
$ws.Range("M116").Value = 1950.3334
$ws.Range("N116").Value = -12364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1551.7954
$ws.Range("I137").Value = 1082.4445
$ws.Range("J137").Value = 1672.4857
$ws.Range("K137").Value = 3247.3335
$ws.Range("L137").Value = 5017.4571
$ws.Range("M137").Value = -697.3335000000002
$ws.Range("N137").Value = -10117.4571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6706.607
$ws.Range("I138").Value = 1940.5625
$ws.Range("J138").Value = 13061.333
$ws.Range("K138").Value = 5821.6875
$ws.Range("L138").Value = 39183.999
$ws.Range("M138").Value = -681.6875
$ws.Range("N138").Value = -49463.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1675.3334
$ws.Range("I61").Value = 781.1786
$ws.Range("J61").Value = 2638.2693
$ws.Range("K61").Value = 781.1786
$ws.Range("L61").Value = 2638.2693
$ws.Range("M61").Value = -569.1786
$ws.Range("N61").Value = -3062.2693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1731.65
$ws.Range("I74").Value = 2371.6
$ws.Range("J74").Value = 1518.3334
$ws.Range("K74").Value = 2371.6
$ws.Range("L74").Value = 1518.3334
$ws.Range("M74").Value = -1497.6
$ws.Range("N74").Value = -3266.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1731.65
$ws.Range("I77").Value = 2371.6
$ws.Range("J77").Value = 1518.3334
$ws.Range("K77").Value = 11858
$ws.Range("L77").Value = 7591.666999999999
$ws.Range("M77").Value = -7490
$ws.Range("N77").Value = -16327.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1675.3334
$ws.Range("I136").Value = 781.1786
$ws.Range("J136").Value = 2638.2693
$ws.Range("K136").Value = 2343.5358
$ws.Range("L136").Value = 7914.8079
$ws.Range("M136").Value = 206.4642000000003
$ws.Range("N136").Value = -13014.8079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1147.5
$ws.Range("I12").Value = 1147.5
$ws.Range("K12").Value = 1147.5
$ws.Range("M12").Value = -979.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 63247.168
$ws.Range("I86").Value = 80275.86
$ws.Range("J86").Value = 3646.75
$ws.Range("K86").Value = 80275.86
$ws.Range("L86").Value = 3646.75
$ws.Range("M86").Value = -79152.86
$ws.Range("N86").Value = -5892.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 63247.168
$ws.Range("I89").Value = 80275.86
$ws.Range("J89").Value = 3646.75
$ws.Range("K89").Value = 401379.3
$ws.Range("L89").Value = 18233.75
$ws.Range("M89").Value = -395763.3
$ws.Range("N89").Value = -29465.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3005
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 870.5454999999999
$ws.Range("I16").Value = 632.4545000000001
$ws.Range("J16").Value = 1108.6364
$ws.Range("K16").Value = 632.4545000000001
$ws.Range("L16").Value = 1108.6364
$ws.Range("M16").Value = -345.4545000000001
$ws.Range("N16").Value = -1682.6364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7066
$ws.Range("I62").Value = 4842
$ws.Range("J62").Value = 9290
$ws.Range("K62").Value = 4842
$ws.Range("L62").Value = 9290
$ws.Range("M62").Value = -4218
$ws.Range("N62").Value = -10538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7066
$ws.Range("I65").Value = 4842
$ws.Range("J65").Value = 9290
$ws.Range("K65").Value = 24210
$ws.Range("L65").Value = 46450
$ws.Range("M65").Value = -21090
$ws.Range("N65").Value = -52690

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2575.375
$ws.Range("I86").Value = 1824.0834
$ws.Range("J86").Value = 3326.6667
$ws.Range("K86").Value = 1824.0834
$ws.Range("L86").Value = 3326.6667
$ws.Range("M86").Value = -701.0834
$ws.Range("N86").Value = -5572.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2575.375
$ws.Range("I89").Value = 1824.0834
$ws.Range("J89").Value = 3326.6667
$ws.Range("K89").Value = 9120.416999999999
$ws.Range("L89").Value = 16633.3335
$ws.Range("M89").Value = -3504.416999999999
$ws.Range("N89").Value = -27865.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1026.45
$ws.Range("I105").Value = 973.6875
$ws.Range("J105").Value = 1237.5
$ws.Range("K105").Value = 973.6875
$ws.Range("L105").Value = 1237.5
$ws.Range("M105").Value = 773.3125
$ws.Range("N105").Value = -4731.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 870.5454999999999
$ws.Range("I113").Value = 632.4545000000001
$ws.Range("J113").Value = 1108.6364
$ws.Range("K113").Value = 632.4545000000001
$ws.Range("L113").Value = 1108.6364
$ws.Range("M113").Value = 1537.5455
$ws.Range("N113").Value = -5448.6364

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 22020.32
$ws.Range("I68").Value = 1938
$ws.Range("J68").Value = 24251.69
$ws.Range("K68").Value = 5814
$ws.Range("L68").Value = 72755.06999999999
$ws.Range("M68").Value = -5003
$ws.Range("N68").Value = -74377.06999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 22020.32
$ws.Range("I71").Value = 1938
$ws.Range("J71").Value = 24251.69
$ws.Range("K71").Value = 17442
$ws.Range("L71").Value = 218265.21
$ws.Range("M71").Value = -13386
$ws.Range("N71").Value = -226377.21

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 26238.25
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 26238.25
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 78714.75
$ws.Range("N80").Value = -80586.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 26238.25
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 26238.25
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 236144.25
$ws.Range("N83").Value = -245504.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 34616492
$ws.Range("I86").Value = 1100
$ws.Range("J86").Value = 56251110
$ws.Range("K86").Value = 3300
$ws.Range("L86").Value = 168753330
$ws.Range("M86").Value = -2114
$ws.Range("N86").Value = -168755702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 34616492
$ws.Range("I89").Value = 1100
$ws.Range("J89").Value = 56251110
$ws.Range("K89").Value = 9900
$ws.Range("L89").Value = 506259990
$ws.Range("M89").Value = -3972
$ws.Range("N89").Value = -506271846

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 86018.30499999999
$ws.Range("I98").Value = 1250
$ws.Range("J98").Value = 101430.73
$ws.Range("K98").Value = 3750
$ws.Range("L98").Value = 304292.19
$ws.Range("M98").Value = -2252
$ws.Range("N98").Value = -307288.19

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 963.0227
$ws.Range("I107").Value = 388.9524
$ws.Range("J107").Value = 1487.174
$ws.Range("K107").Value = 1166.8572
$ws.Range("L107").Value = 4461.522
$ws.Range("M107").Value = 753.1428000000001
$ws.Range("N107").Value = -8301.522000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 797.67
$ws.Range("I131").Value = 315
$ws.Range("J131").Value = 823.07367
$ws.Range("K131").Value = 945
$ws.Range("L131").Value = 2469.22101
$ws.Range("M131").Value = 4095
$ws.Range("N131").Value = -12549.22101

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 109927.21
$ws.Range("I70").Value = 171088.92
$ws.Range("J70").Value = 5078.5713
$ws.Range("K70").Value = 171088.92
$ws.Range("L70").Value = 5078.5713
$ws.Range("M70").Value = -170818.92
$ws.Range("N70").Value = -5618.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 109927.21
$ws.Range("I73").Value = 171088.92
$ws.Range("J73").Value = 5078.5713
$ws.Range("K73").Value = 171088.92
$ws.Range("L73").Value = 5078.5713
$ws.Range("M73").Value = -170152.92
$ws.Range("N73").Value = -6950.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 235002.5
$ws.Range("I14").Value = 400000
$ws.Range("J14").Value = 70005
$ws.Range("K14").Value = 400000
$ws.Range("L14").Value = 70005
$ws.Range("M14").Value = -399828
$ws.Range("N14").Value = -70349

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 915.63635
$ws.Range("I93").Value = 915.63635
$ws.Range("K93").Value = 915.63635
$ws.Range("M93").Value = 332.36365

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 9240
$ws.Range("J25").Value = 9240
$ws.Range("L25").Value = 9240
$ws.Range("N25").Value = -9826

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3660
$ws.Range("I126").Value = 3990
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 11970
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -9500
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2545.2273
$ws.Range("I132").Value = 2644.4644
$ws.Range("K132").Value = 7933.3932
$ws.Range("M132").Value = -5403.3932
